$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test rows (20-28): three groups of three runs each, mirroring the
# existing "3-run group" pattern already present in the log. Enter
# column B (group) first for every row, then column A (name) for every
# row -- matches the order the new shared strings show up in.
$ws.Range("B20").Value = "15.Start"
$ws.Range("B21").Value = "15.Start"
$ws.Range("B22").Value = "15.Start"
$ws.Range("B23").Value = "15.Minute"
$ws.Range("B24").Value = "15.Minute"
$ws.Range("B25").Value = "15.Minute"
$ws.Range("B26").Value = "15.Mix"
$ws.Range("B27").Value = "15.Mix"
$ws.Range("B28").Value = "15.Mix"

$ws.Range("A20").Value = "15ST.1"
$ws.Range("A21").Value = "15ST.2"
$ws.Range("A22").Value = "15ST.3"
$ws.Range("A23").Value = "15MT.1"
$ws.Range("A24").Value = "15MT.2"
$ws.Range("A25").Value = "15MT.3"
$ws.Range("A26").Value = "15MX.1"
$ws.Range("A27").Value = "15MX.2"
$ws.Range("A28").Value = "15MX.3"

# D/E/F carry the same formulas used throughout the sheet (random
# chance, change, group change) -- fill them in one row at a time so
# each keeps the column's existing number format.
for ($row = 20; $row -le 28; $row++) {
    $ws.Range("D$row").Formula = "=1/3"
    $ws.Range("D$row").NumberFormat = "0.00%"
    $ws.Range("E$row").Formula = "=C$row/D$row"
    $ws.Range("E$row").NumberFormat = "0.00"
    $ws.Range("F$row").Formula = "=AVERAGEIF(B:B,B$row,E:E)"
    $ws.Range("F$row").NumberFormat = "0.00"
}
$ws.Range("G24").Select()
